# Apply the "Jogos da Semana" update:
#  - Row 2's match data is replaced with a new fixture.
#  - Columns BC:BD (Odd_CS_3-3_HT / Odd_CS_4-4_HT) are removed entirely.
#  - Row 3 (the Brazil Serie A fixture) is removed entirely.
#  - Dimension shrinks from A1:BD3 to A1:BB2 as a consequence.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove columns BC:BD (entire columns, so nothing else shifts unexpectedly).
$ws.Range("BC1:BD1").EntireColumn.Delete()

# 2) Remove row 3 entirely.
$ws.Rows(3).Delete()

# 3) Update row 2 with the new match's data.
$ws.Range("A2").Value = "CEYAc3cC"
$ws.Range("C2").Value = "20:30"
$ws.Range("D2").Value = "VENEZUELA - LIGA FUTVE"
$ws.Range("E2").Value = "Estudiantes Merida"
$ws.Range("F2").Value = "La Guaira"
$ws.Range("G2").Value = 2.37
$ws.Range("H2").Value = 2.95
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 2.9
$ws.Range("K2").Value = 2.02
$ws.Range("L2").Value = 3.55
$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 8
$ws.Range("P2").Value = 2.77
$ws.Range("Q2").Value = 1.98
$ws.Range("R2").Value = 1.65
$ws.Range("S2").Value = 1.47
$ws.Range("T2").Value = 2.35
$ws.Range("U2").Value = 1.7
$ws.Range("V2").Value = 1.91
$ws.Range("W2").Value = 7.7
$ws.Range("X2").Value = 12
$ws.Range("Z2").Value = 26
$ws.Range("AA2").Value = 19.5
$ws.Range("AC2").Value = 8.25
$ws.Range("AD2").Value = 5.8
$ws.Range("AE2").Value = 13
$ws.Range("AF2").Value = 65
$ws.Range("AG2").Value = 500
$ws.Range("AH2").Value = 8.5
$ws.Range("AI2").Value = 15.5
$ws.Range("AJ2").Value = 10.75
$ws.Range("AK2").Value = 40
$ws.Range("AL2").Value = 28
$ws.Range("AM2").Value = 35
$ws.Range("AN2").Value = 4.3
$ws.Range("AP2").Value = 18
$ws.Range("AQ2").Value = 50
$ws.Range("AR2").Value = 75
$ws.Range("AS2").Value = 200
$ws.Range("AT2").Value = 2.52
$ws.Range("AU2").Value = 6.5
$ws.Range("AV2").Value = 50
$ws.Range("AW2").Value = 5
$ws.Range("AX2").Value = 16.5
$ws.Range("AY2").Value = 22
$ws.Range("AZ2").Value = 80
$ws.Range("BA2").Value = 110
$ws.Range("BB2").Value = 250
